# Rename the worksheet tab to reflect the state abbreviation (Maine = ME)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ME"
